$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 130.83333
$ws.Range("I9").Value = 139.25
$ws.Range("J9").Value = 114
$ws.Range("K9").Value = 139.25
$ws.Range("L9").Value = 114
$ws.Range("M9").Value = 29.75
$ws.Range("N9").Value = -452
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H86").Value = 4552.2
$ws.Range("I86").Value = 3332.3333
$ws.Range("K86").Value = 3332.3333
$ws.Range("M86").Value = -2209.3333
$ws.Range("H89").Value = 4552.2
$ws.Range("I89").Value = 3332.3333
$ws.Range("K89").Value = 16661.6665
$ws.Range("M89").Value = -11045.6665
$ws.Range("H92").Value = 791.7143
$ws.Range("I92").Value = 698.7692
$ws.Range("K92").Value = 698.7692
$ws.Range("M92").Value = 549.2308
$ws.Range("H96").Value = 2928.7778
$ws.Range("J96").Value = 3689.6
$ws.Range("L96").Value = 11068.8
$ws.Range("N96").Value = -13814.8
$ws.Range("H111").Value = 2733
$ws.Range("I111").Value = 3013.6
$ws.Range("K111").Value = 9040.799999999999
$ws.Range("M111").Value = -5973.799999999999
$ws.Range("H129").Value = 200000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 200000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 600000
$ws.Range("N129").Value = -610000
$ws.Range("H141").Value = 1776.1333
$ws.Range("I141").Value = 1045.8572
$ws.Range("K141").Value = 3137.5716
$ws.Range("M141").Value = 2042.4284
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("M129").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 982.1667
$ws.Range("I2").Value = 1495.3334
$ws.Range("J2").Value = 469
$ws.Range("K2").Value = 1495.3334
$ws.Range("L2").Value = 469
$ws.Range("M2").Value = -1382.3334
$ws.Range("N2").Value = -695
$ws.Range("H88").Value = 412.0909
$ws.Range("I88").Value = 412.0909
$ws.Range("K88").Value = 412.0909
$ws.Range("M88").Value = -6.090899999999976
$ws.Range("H91").Value = 412.0909
$ws.Range("I91").Value = 412.0909
$ws.Range("K91").Value = 412.0909
$ws.Range("M91").Value = 991.9091000000001
$ws.Range("H97").Value = 1501.25
$ws.Range("I97").Value = 1503.3334
$ws.Range("J97").Value = 1495
$ws.Range("K97").Value = 1503.3334
$ws.Range("L97").Value = 1495
$ws.Range("M97").Value = -1007.3334
$ws.Range("N97").Value = -2487
$ws.Range("H102").Value = 806.6
$ws.Range("I102").Value = 507.6154
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 507.6154
$ws.Range("L102").Value = 2750
$ws.Range("M102").Value = 1114.3846
$ws.Range("N102").Value = -5994
$ws.Range("H116").Value = 982.1667
$ws.Range("I116").Value = 1495.3334
$ws.Range("J116").Value = 469
$ws.Range("K116").Value = 1495.3334
$ws.Range("L116").Value = 469
$ws.Range("M116").Value = 798.6666
$ws.Range("N116").Value = -5057
$ws.Range("H132").Value = 1728.15
$ws.Range("I132").Value = 1650.7106
$ws.Range("K132").Value = 4952.1318
$ws.Range("M132").Value = -2422.1318

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 982.1667
$ws.Range("I3").Value = 1495.3334
$ws.Range("J3").Value = 469
$ws.Range("K3").Value = 1495.3334
$ws.Range("L3").Value = 469
$ws.Range("M3").Value = -1381.3334
$ws.Range("N3").Value = -697
$ws.Range("H134").Value = 1424.9
$ws.Range("I134").Value = 1236.7894
$ws.Range("K134").Value = 3710.3682
$ws.Range("M134").Value = -1175.3682

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1203.9286
$ws.Range("I16").Value = 690.7
$ws.Range("K16").Value = 690.7
$ws.Range("M16").Value = -403.7
$ws.Range("H103").Value = 11485.667
$ws.Range("I103").Value = 12982.8
$ws.Range("J103").Value = 4000
$ws.Range("K103").Value = 12982.8
$ws.Range("L103").Value = 4000
$ws.Range("M103").Value = -11810.8
$ws.Range("N103").Value = -6344
$ws.Range("H107").Value = 747.9524
$ws.Range("I107").Value = 536.6667
$ws.Range("K107").Value = 536.6667
$ws.Range("M107").Value = 1383.3333
$ws.Range("H113").Value = 1203.9286
$ws.Range("I113").Value = 690.7
$ws.Range("K113").Value = 690.7
$ws.Range("M113").Value = 1479.3
$ws.Range("H132").Value = 2679.524
$ws.Range("I132").Value = 1758.5714
$ws.Range("J132").Value = 4521.4287
$ws.Range("K132").Value = 5275.7142
$ws.Range("L132").Value = 13564.2861
$ws.Range("M132").Value = -2745.7142
$ws.Range("N132").Value = -18624.2861
$ws.Range("H134").Value = 3440
$ws.Range("I134").Value = 2692
$ws.Range("K134").Value = 8076
$ws.Range("M134").Value = -5541

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 538.3333
$ws.Range("H70").Value = 300
$ws.Range("I70").Value = 300
$ws.Range("K70").Value = 900
$ws.Range("M70").Value = -585
$ws.Range("H73").Value = 300
$ws.Range("I73").Value = 300
$ws.Range("K73").Value = 900
$ws.Range("M73").Value = 192
$ws.Range("H132").Value = 1824.8182
$ws.Range("I132").Value = 1218.6
$ws.Range("K132").Value = 10967.4
$ws.Range("M132").Value = -8437.4
$ws.Range("H137").Value = 7827.533
$ws.Range("I137").Value = 5999.4
$ws.Range("J137").Value = 8741.6
$ws.Range("K137").Value = 17998.2
$ws.Range("L137").Value = 26224.8
$ws.Range("M137").Value = -12898.2
$ws.Range("N137").Value = -36424.8
$ws.Range("H139").Value = 4415.278
$ws.Range("I139").Value = 2089.182
$ws.Range("J139").Value = 8070.5713
$ws.Range("K139").Value = 6267.545999999999
$ws.Range("L139").Value = 24211.7139
$ws.Range("M139").Value = -1127.545999999999
$ws.Range("N139").Value = -34491.7139

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5974.375
$ws.Range("J36").Value = 4669.4
$ws.Range("L36").Value = 4669.4
$ws.Range("N36").Value = -5639.4
$ws.Range("H80").Value = 14248.25
$ws.Range("I80").Value = 13499
$ws.Range("J80").Value = 14997.5
$ws.Range("K80").Value = 13499
$ws.Range("L80").Value = 14997.5
$ws.Range("M80").Value = -12501
$ws.Range("N80").Value = -16993.5
$ws.Range("H83").Value = 14248.25
$ws.Range("I83").Value = 13499
$ws.Range("J83").Value = 14997.5
$ws.Range("K83").Value = 67495
$ws.Range("L83").Value = 74987.5
$ws.Range("M83").Value = -62503
$ws.Range("N83").Value = -84971.5
$ws.Range("H102").Value = 1728.5555
$ws.Range("I102").Value = 470.30768
$ws.Range("K102").Value = 470.30768
$ws.Range("M102").Value = 1151.69232
$ws.Range("H107").Value = 903.1852
$ws.Range("I107").Value = 980.5714
$ws.Range("J107").Value = 819.8461
$ws.Range("K107").Value = 980.5714
$ws.Range("L107").Value = 819.8461
$ws.Range("M107").Value = 939.4286
$ws.Range("N107").Value = -4659.8461
$ws.Range("H132").Value = 1986.25
$ws.Range("I132").Value = 1526.1333
$ws.Range("K132").Value = 4578.3999
$ws.Range("M132").Value = -2048.3999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6173.75
$ws.Range("I61").Value = 7380.3335
$ws.Range("J61").Value = 4967.1665
$ws.Range("K61").Value = 7380.3335
$ws.Range("L61").Value = 4967.1665
$ws.Range("M61").Value = -7178.3335
$ws.Range("N61").Value = -5371.1665
$ws.Range("H68").Value = 3207.75
$ws.Range("I68").Value = 2666.6667
$ws.Range("J68").Value = 4831
$ws.Range("K68").Value = 2666.6667
$ws.Range("L68").Value = 4831
$ws.Range("M68").Value = -1917.6667
$ws.Range("N68").Value = -6329
$ws.Range("H71").Value = 3207.75
$ws.Range("I71").Value = 2666.6667
$ws.Range("J71").Value = 4831
$ws.Range("K71").Value = 13333.3335
$ws.Range("L71").Value = 24155
$ws.Range("M71").Value = -9589.333500000001
$ws.Range("N71").Value = -31643
$ws.Range("H113").Value = 6173.75
$ws.Range("I113").Value = 7380.3335
$ws.Range("J113").Value = 4967.1665
$ws.Range("K113").Value = 7380.3335
$ws.Range("L113").Value = 4967.1665
$ws.Range("M113").Value = -5210.3335
$ws.Range("N113").Value = -9307.166499999999
$ws.Range("H122").Value = 9926.117
$ws.Range("I122").Value = 9609
$ws.Range("K122").Value = 28827
$ws.Range("M122").Value = -26377
$ws.Range("H132").Value = 4089.2666
$ws.Range("I132").Value = 3449.0454
$ws.Range("K132").Value = 10347.1362
$ws.Range("M132").Value = -7817.136200000001
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 2100
$ws.Range("I52").Value = 2100
$ws.Range("K52").Value = 2100
$ws.Range("M52").Value = -1874
$ws.Range("H62").Value = 7861.067
$ws.Range("I62").Value = 6000
$ws.Range("K62").Value = 6000
$ws.Range("M62").Value = -5376
$ws.Range("H65").Value = 7861.067
$ws.Range("I65").Value = 6000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26880
